$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.414.68'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.68'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9987'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6331'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.0000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07577'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2966'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07711'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.865.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.005'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6868'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001006'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.129.62'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.167'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.446.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.582'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1400'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.391'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.469'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05716'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.51%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.259'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.129'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.029'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.851'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.157'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7175'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.589'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.265.60'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01812'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.781'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9072'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.187'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.29'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.16%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.092'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4055'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000118'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.129'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.688'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.78%  '
